$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu results for the 380 kV case (rows 2-25, columns B-F and I-N)
$data = @{
    2 = @(1.02,1.030207184111771,1.034257598597455,0.9926147277508489,1.040815895542663,1.037413502484381,1.035349705075748,1.037057591274386,0.9955398523336033,1.043597172926108,1.036820020413245)
    3 = @(1.02,1.031031465429633,1.034867607726591,0.9936372048519304,1.041910106369248,1.037634021521089,1.035815982329411,1.037477281658688,0.9963617723202692,1.044501138628534,1.037286959834073)
    4 = @(1.02,1.031565005212316,1.035262410011656,0.9942998659930995,1.042618681233364,1.037775440565625,1.036117210774419,1.03774824143123,0.9968940712668345,1.045086022395363,1.037588616058024)
    5 = @(1.02,1.031789345921192,1.035428403961425,0.9945786998346017,1.042916696524642,1.037834588123642,1.036243730397166,1.037862006422149,0.997117960005301,1.045331896711596,1.037715315353147)
    6 = @(1.02,1.031827016057556,1.035456276152835,0.9946255319796338,1.042966742233701,1.037844501355498,1.036264966721743,1.037881099446435,0.9971555583673453,1.045373179437988,1.037736581835741)
    7 = @(1.02,1.031568002705916,1.0352646279578,0.9943035907982488,1.042622662816566,1.037776232097419,1.036118901795966,1.037749762141432,0.9968970624462087,1.045089307823329,1.037590309481015)
    8 = @(1.02,1.030485716941353,1.034463735463888,0.9929600610674301,1.041185575529866,1.037488290838033,1.035507385345437,1.037199552919305,0.995817528259106,1.043902680775195,1.036977924606984)
    9 = @(1.02,1.028579982860312,1.033053178634189,0.9906006454969559,1.038657462155874,1.036971189855706,1.034426147823158,1.036225399383119,0.9939188001724441,1.041811395250575,1.035895151603403)
    10 = @(1.02,1.027310501035343,1.03211337402536,0.989033133672735,1.036974922589244,1.036619961652059,1.033702910765472,1.035572916773834,0.9926553831429383,1.040417041492662,1.03517088746614)
    11 = @(1.02,1.026761055174651,1.031706580160021,0.988355674866747,1.036247051251263,1.036466342903341,1.033389180129922,1.035289673454734,0.9921088820399291,1.039813239188055,1.034856711297101)
    12 = @(1.02,1.026557004949189,1.031555502395082,0.9881042295826724,1.035976789918391,1.036409052051524,1.033272562697116,1.035184357897538,0.9919059725120875,1.039588954914782,1.034739928254165)
    13 = @(1.02,1.026600772657654,1.031587907998713,0.9881581567098651,1.036034757247043,1.036421351533253,1.033297581315424,1.035206953241416,0.9919494934313052,1.039637064880836,1.034764982401781)
    14 = @(1.02,1.026744187528722,1.031694091535585,0.9883348863814464,1.036224709271968,1.03646161191258,1.033379542203225,1.035280970202785,0.9920921077337197,1.03979469988359,1.034847059683444)
    15 = @(1.02,1.026832555291001,1.031759517884565,0.9884438009545853,1.036341758577318,1.036486387199137,1.033430029918649,1.035326560426447,0.9921799884222134,1.03989182333518,1.034897619097214)
    16 = @(1.02,1.027346971273886,1.032140374826103,0.9890781214508737,1.037023243362537,1.036630124536803,1.033723720220652,1.03559169973286,0.9926916645766087,1.040457113105307,1.035191726473133)
    17 = @(1.02,1.02766971807702,1.032379316861052,0.989476357848556,1.037450902952287,1.036719876698706,1.033907793939462,1.035757824081012,0.9930127773699352,1.040811694630994,1.035376061597738)
    18 = @(1.02,1.0278579946923,1.032518701868228,0.9897087662937556,1.037700415079293,1.036772079510869,1.034015106528929,1.035854652647599,0.9932001317071769,1.041018512343894,1.035483526583391)
    19 = @(1.02,1.027922196165614,1.032566230913075,0.9897880325774034,1.037785503367905,1.03678985418279,1.034051688078967,1.035887656963145,0.9932640239640975,1.041089031198408,1.035520160083425)
    20 = @(1.02,1.027635087922634,1.032353679167225,0.9894336180360679,1.037405012345098,1.036710262456383,1.033888050202068,1.035740007636469,0.9929783193494215,1.040773651767316,1.035356289821972)
    21 = @(1.02,1.026701954353586,1.031662822468183,0.9882828385668249,1.036168770302436,1.036449762579275,1.033355409073892,1.035259176981796,0.9920501090198102,1.039748280444225,1.034822892282299)
    22 = @(1.02,1.026115479163878,1.031228590607866,0.9875604150241495,1.035392089052654,1.03628464548091,1.033020031428002,1.034956244644388,0.9914670000341481,1.039103558877795,1.034487038361683)
    23 = @(1.02,1.026426359243969,1.031458771666838,0.9879432794643023,1.035803766110218,1.036372303109703,1.033197867151413,1.035116892804942,0.991776070289318,1.039445340651954,1.034665126632219)
    24 = @(1.02,1.02765073572512,1.032365263693903,0.9894529299347244,1.037425748132627,1.036714607176311,1.033896971715598,1.035748058334003,0.9929938892766442,1.040790841711242,1.035365224005074)
    25 = @(1.02,1.029072488470655,1.03341774763663,0.9912096547607049,1.039310536773941,1.037106020094965,1.034706102705892,1.036477782531281,0.9944092447426414,1.042352074155126,1.036175504054183)
}

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

foreach ($rowKey in $data.Keys) {
    $rowValues = $data[$rowKey]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $rowKey).Value = $rowValues[$i]
    }
}
